# Word COM-interop script: mark a batch of "Workflow document" bullet
# points as done by colouring them with the same accent6/BF green that
# is already used elsewhere in the document
# (w:color w:val="538135" w:themeColor="accent6" w:themeShade="BF").
#
# Strategy: for each target paragraph we pull its OOXML fragment via
# Range.WordOpenXML, splice in the <w:rPr><w:color .../></w:rPr> using
# regex (so the existing w14:paraId / rsid / pPr content - pStyle,
# numPr, etc. - is preserved byte for byte), then write it back with
# Range.InsertXML on the very same Range (this replaces only that
# paragraph's contents and does not disturb neighbouring paragraphs).

$d = $word.ActiveDocument

$colorRpr = '<w:rPr><w:color w:val="538135" w:themeColor="accent6" w:themeShade="BF"/></w:rPr>'

function Get-ParaFragment($range) {
    $xml = $range.WordOpenXML
    $m = [regex]::Match($xml, '<w:p[ >][\s\S]*?</w:p>')
    return $m.Value
}

function Set-XmlFragment($range, $frag) {
    $wrapper = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
        '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
        '<pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:body>' +
        $frag + '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
    $null = $range.InsertXML($wrapper)
}

# Colours every run in the paragraph (and the paragraph mark's rPr,
# inside pPr) with the standard green. Used for paragraphs whose text
# is entirely one colour.
function Set-ParagraphGreen($idx) {
    $p = $d.Paragraphs.Item($idx)
    $full = $p.Range
    $frag = Get-ParaFragment $full

    # Paragraph-mark formatting: add <w:rPr>...</w:rPr> as the last
    # child of <w:pPr> (pPr only ever closes with </w:pPr> once).
    $frag = $frag -replace '</w:pPr>', ($colorRpr + '</w:pPr>')

    # Colour every run-opening tag, regardless of any rsid attributes
    # Word's object model may have put on it.
    $frag = [regex]::Replace($frag, '<w:r( [^>]*)?>', ('<w:r>' + $colorRpr))

    Set-XmlFragment $full $frag
}

35, 36, 37, 45, 46, 47, 52, 53, 54, 56, 57, 58, 59 | ForEach-Object {
    Set-ParagraphGreen $_
}

# Paragraph 55 ("Game stage that will hold multiple waves, base reward
# money on complete, star ranking") is special: its pPr stays
# untouched, its first two logical runs turn green but the trailing
# "star ranking" stays black, which means the second original run has
# to be split in two.
$p55 = $d.Paragraphs.Item(55)
$full55 = $p55.Range
$xml55 = $full55.WordOpenXML
$popen = [regex]::Match($xml55, '<w:p[^>]*>').Value
$pPrMatch = [regex]::Match($xml55, '<w:pPr>[\s\S]*?</w:pPr>')
$pPr = $pPrMatch.Value

$run1 = '<w:r>' + $colorRpr + '<w:t>Game stage that will hold multiple waves</w:t></w:r>'
$run2 = '<w:r>' + $colorRpr + '<w:t xml:space="preserve">, base reward money on complete, </w:t></w:r>'
$run3 = '<w:r><w:t>star ranking</w:t></w:r>'

$newFrag55 = $popen + $pPr + $run1 + $run2 + $run3 + '</w:p>'
Set-XmlFragment $full55 $newFrag55

Write-Output "done"
